$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 86: G1 / Test1
$ws.Cells.Item(86, 1).Value = "G1"
$ws.Cells.Item(86, 2).Value = "Test1"
$ws.Cells.Item(86, 3).Value = 45903
$ws.Cells.Item(86, 3).NumberFormat = $ws.Cells.Item(85, 3).NumberFormat
$ws.Cells.Item(86, 4).Value = 0.665003107782612
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 6).Value = -0.01

# Row 87: G2 / sedrftgyhuioygtfrd
$ws.Cells.Item(87, 1).Value = "G2"
$ws.Cells.Item(87, 2).Value = "sedrftgyhuioygtfrd"
$ws.Cells.Item(87, 3).Value = 45903
$ws.Cells.Item(87, 3).NumberFormat = $ws.Cells.Item(85, 3).NumberFormat
$ws.Cells.Item(87, 4).Value = 0.665003107782612
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = -0.01
